$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9,8).Value = 200
$ws.Cells.Item(9,9).Value = 0
$ws.Cells.Item(9,10).Value = 200
$ws.Cells.Item(9,11).Value = 0
$ws.Cells.Item(9,12).Value = 200
$ws.Cells.Item(9,13).Value = ""
$ws.Cells.Item(9,14).Value = -538
$ws.Cells.Item(18,8).Value = 887
$ws.Cells.Item(18,9).Value = 887
$ws.Cells.Item(18,11).Value = 887
$ws.Cells.Item(18,13).Value = -603
$ws.Cells.Item(31,8).Value = 192.66667
$ws.Cells.Item(31,9).Value = 192.66667
$ws.Cells.Item(31,10).Value = 0
$ws.Cells.Item(31,11).Value = 578.00001
$ws.Cells.Item(31,12).Value = 0
$ws.Cells.Item(31,13).Value = -348.00001
$ws.Cells.Item(31,14).Value = ""
$ws.Cells.Item(62,8).Value = 4433.143
$ws.Cells.Item(62,9).Value = 4198.3335
$ws.Cells.Item(62,10).Value = 4609.25
$ws.Cells.Item(62,11).Value = 4198.3335
$ws.Cells.Item(62,12).Value = 4609.25
$ws.Cells.Item(62,13).Value = -3574.3335
$ws.Cells.Item(62,14).Value = -5857.25
$ws.Cells.Item(65,8).Value = 4433.143
$ws.Cells.Item(65,9).Value = 4198.3335
$ws.Cells.Item(65,10).Value = 4609.25
$ws.Cells.Item(65,11).Value = 20991.6675
$ws.Cells.Item(65,12).Value = 23046.25
$ws.Cells.Item(65,13).Value = -17871.6675
$ws.Cells.Item(65,14).Value = -29286.25
$ws.Cells.Item(86,8).Value = 3843
$ws.Cells.Item(86,9).Value = 5141.3335
$ws.Cells.Item(86,10).Value = 1895.5
$ws.Cells.Item(86,11).Value = 5141.3335
$ws.Cells.Item(86,12).Value = 1895.5
$ws.Cells.Item(86,13).Value = -4018.3335
$ws.Cells.Item(86,14).Value = -4141.5
$ws.Cells.Item(88,8).Value = 900
$ws.Cells.Item(88,10).Value = 600
$ws.Cells.Item(88,12).Value = 600
$ws.Cells.Item(88,14).Value = -1412
$ws.Cells.Item(89,8).Value = 3843
$ws.Cells.Item(89,9).Value = 5141.3335
$ws.Cells.Item(89,10).Value = 1895.5
$ws.Cells.Item(89,11).Value = 25706.6675
$ws.Cells.Item(89,12).Value = 9477.5
$ws.Cells.Item(89,13).Value = -20090.6675
$ws.Cells.Item(89,14).Value = -20709.5
$ws.Cells.Item(91,8).Value = 900
$ws.Cells.Item(91,10).Value = 600
$ws.Cells.Item(91,12).Value = 600
$ws.Cells.Item(91,14).Value = -3408
$ws.Cells.Item(137,8).Value = 2901.1667
$ws.Cells.Item(137,10).Value = 3281.4
$ws.Cells.Item(137,12).Value = 9844.200000000001
$ws.Cells.Item(137,14).Value = -14944.2
$ws.Cells.Item(138,8).Value = 3621.7
$ws.Cells.Item(138,10).Value = 3621.7
$ws.Cells.Item(138,12).Value = 10865.1
$ws.Cells.Item(138,14).Value = -21145.1
$ws.Cells.Item(141,8).Value = 2249
$ws.Cells.Item(141,9).Value = 2249
$ws.Cells.Item(141,11).Value = 6747
$ws.Cells.Item(141,13).Value = -1567

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45,8).Value = 2020.1818
$ws.Cells.Item(45,9).Value = 1746.1428
$ws.Cells.Item(45,11).Value = 1746.1428
$ws.Cells.Item(45,13).Value = -1369.1428
$ws.Cells.Item(61,8).Value = 5596.6
$ws.Cells.Item(61,9).Value = 5596.6
$ws.Cells.Item(61,11).Value = 5596.6
$ws.Cells.Item(61,13).Value = -5384.6
$ws.Cells.Item(74,8).Value = 3621.7778
$ws.Cells.Item(74,9).Value = 3621.7778
$ws.Cells.Item(74,11).Value = 3621.7778
$ws.Cells.Item(74,13).Value = -2747.7778
$ws.Cells.Item(77,8).Value = 3621.7778
$ws.Cells.Item(77,9).Value = 3621.7778
$ws.Cells.Item(77,11).Value = 18108.889
$ws.Cells.Item(77,13).Value = -13740.889
$ws.Cells.Item(102,8).Value = 1780
$ws.Cells.Item(102,9).Value = 1510
$ws.Cells.Item(102,11).Value = 1510
$ws.Cells.Item(102,13).Value = 112
$ws.Cells.Item(132,8).Value = 1781.8823
$ws.Cells.Item(132,9).Value = 1714.8462
$ws.Cells.Item(132,11).Value = 5144.5386
$ws.Cells.Item(132,13).Value = -2614.5386
$ws.Cells.Item(136,8).Value = 5596.6
$ws.Cells.Item(136,9).Value = 5596.6
$ws.Cells.Item(136,11).Value = 16789.8
$ws.Cells.Item(136,13).Value = -14239.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105,8).Value = 2626.182
$ws.Cells.Item(105,9).Value = 1598.6666
$ws.Cells.Item(105,10).Value = 7250
$ws.Cells.Item(105,11).Value = 1598.6666
$ws.Cells.Item(105,12).Value = 7250
$ws.Cells.Item(105,13).Value = 148.3334
$ws.Cells.Item(105,14).Value = -10744
$ws.Cells.Item(134,8).Value = 5058.857
$ws.Cells.Item(134,9).Value = 5058.857
$ws.Cells.Item(134,11).Value = 15176.571
$ws.Cells.Item(134,13).Value = -12641.571

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 1963.3726
$ws.Cells.Item(31,9).Value = 1500.2632
$ws.Cells.Item(31,10).Value = 3317.077
$ws.Cells.Item(31,11).Value = 1500.2632
$ws.Cells.Item(31,12).Value = 3317.077
$ws.Cells.Item(31,13).Value = -1205.2632
$ws.Cells.Item(31,14).Value = -3907.077
$ws.Cells.Item(34,8).Value = 1963.3726
$ws.Cells.Item(34,9).Value = 1500.2632
$ws.Cells.Item(34,10).Value = 3317.077
$ws.Cells.Item(34,11).Value = 1500.2632
$ws.Cells.Item(34,12).Value = 3317.077
$ws.Cells.Item(34,13).Value = -1298.2632
$ws.Cells.Item(34,14).Value = -3721.077
$ws.Cells.Item(58,8).Value = 5856.231
$ws.Cells.Item(58,9).Value = 5121.091
$ws.Cells.Item(58,10).Value = 9899.5
$ws.Cells.Item(58,11).Value = 5121.091
$ws.Cells.Item(58,12).Value = 9899.5
$ws.Cells.Item(58,13).Value = -4918.091
$ws.Cells.Item(58,14).Value = -10305.5
$ws.Cells.Item(86,8).Value = 9045.154
$ws.Cells.Item(86,9).Value = 8269.9
$ws.Cells.Item(86,11).Value = 8269.9
$ws.Cells.Item(86,13).Value = -7146.9
$ws.Cells.Item(89,8).Value = 9045.154
$ws.Cells.Item(89,9).Value = 8269.9
$ws.Cells.Item(89,11).Value = 41349.5
$ws.Cells.Item(89,13).Value = -35733.5
$ws.Cells.Item(122,8).Value = 1078
$ws.Cells.Item(122,9).Value = 937.3333
$ws.Cells.Item(122,10).Value = 1500
$ws.Cells.Item(122,11).Value = 2811.9999
$ws.Cells.Item(122,12).Value = 4500
$ws.Cells.Item(122,13).Value = -361.9998999999998
$ws.Cells.Item(122,14).Value = -9400
$ws.Cells.Item(132,8).Value = 2005.7142
$ws.Cells.Item(132,9).Value = 1756.75
$ws.Cells.Item(132,11).Value = 5270.25
$ws.Cells.Item(132,13).Value = -2740.25
$ws.Cells.Item(134,8).Value = 103162
$ws.Cells.Item(134,9).Value = 167937
$ws.Cells.Item(134,10).Value = 5999.5
$ws.Cells.Item(134,11).Value = 503811
$ws.Cells.Item(134,12).Value = 17998.5
$ws.Cells.Item(134,13).Value = -501276
$ws.Cells.Item(134,14).Value = -23068.5
$ws.Cells.Item(136,8).Value = 5856.231
$ws.Cells.Item(136,9).Value = 5121.091
$ws.Cells.Item(136,10).Value = 9899.5
$ws.Cells.Item(136,11).Value = 15363.273
$ws.Cells.Item(136,12).Value = 29698.5
$ws.Cells.Item(136,13).Value = -12813.273
$ws.Cells.Item(136,14).Value = -34798.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68,8).Value = 1457.3334
$ws.Cells.Item(68,10).Value = 1666.6666
$ws.Cells.Item(68,12).Value = 4999.9998
$ws.Cells.Item(68,14).Value = -6621.9998
$ws.Cells.Item(71,8).Value = 1457.3334
$ws.Cells.Item(71,10).Value = 1666.6666
$ws.Cells.Item(71,12).Value = 14999.9994
$ws.Cells.Item(71,14).Value = -23111.9994
$ws.Cells.Item(114,8).Value = 445.66666
$ws.Cells.Item(114,9).Value = 445.66666
$ws.Cells.Item(114,10).Value = 0
$ws.Cells.Item(114,11).Value = 1336.99998
$ws.Cells.Item(114,12).Value = 0
$ws.Cells.Item(114,13).Value = 1917.00002
$ws.Cells.Item(114,14).Value = ""

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(57,8).Value = 22600
$ws.Cells.Item(57,10).Value = 25250
$ws.Cells.Item(57,12).Value = 25250
$ws.Cells.Item(57,14).Value = -26890
$ws.Cells.Item(80,8).Value = 3431.6
$ws.Cells.Item(80,10).Value = 4466.5
$ws.Cells.Item(80,12).Value = 4466.5
$ws.Cells.Item(80,14).Value = -6462.5
$ws.Cells.Item(83,8).Value = 3431.6
$ws.Cells.Item(83,10).Value = 4466.5
$ws.Cells.Item(83,12).Value = 22332.5
$ws.Cells.Item(83,14).Value = -32316.5
$ws.Cells.Item(103,8).Value = 0
$ws.Cells.Item(103,10).Value = 0
$ws.Cells.Item(103,12).Value = 0
$ws.Cells.Item(103,14).Value = ""
$ws.Cells.Item(126,8).Value = 9000
$ws.Cells.Item(126,9).Value = 9000
$ws.Cells.Item(126,11).Value = 27000
$ws.Cells.Item(126,13).Value = -24530
$ws.Cells.Item(132,8).Value = 2198.6
$ws.Cells.Item(132,9).Value = 2244.75
$ws.Cells.Item(132,11).Value = 6734.25
$ws.Cells.Item(132,13).Value = -4204.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7,8).Value = 2833.3333
$ws.Cells.Item(7,9).Value = 3000
$ws.Cells.Item(7,11).Value = 3000
$ws.Cells.Item(7,13).Value = -2888
$ws.Cells.Item(16,8).Value = 2886.7778
$ws.Cells.Item(16,9).Value = 3163.3333
$ws.Cells.Item(16,10).Value = 2333.6667
$ws.Cells.Item(16,11).Value = 3163.3333
$ws.Cells.Item(16,12).Value = 2333.6667
$ws.Cells.Item(16,13).Value = -2993.3333
$ws.Cells.Item(16,14).Value = -2673.6667
$ws.Cells.Item(68,8).Value = 26200.2
$ws.Cells.Item(68,10).Value = 39667.668
$ws.Cells.Item(68,12).Value = 39667.668
$ws.Cells.Item(68,14).Value = -41165.668
$ws.Cells.Item(71,8).Value = 26200.2
$ws.Cells.Item(71,10).Value = 39667.668
$ws.Cells.Item(71,12).Value = 198338.34
$ws.Cells.Item(71,14).Value = -205826.34
$ws.Cells.Item(100,8).Value = 3400
$ws.Cells.Item(100,9).Value = 3400
$ws.Cells.Item(100,11).Value = 3400
$ws.Cells.Item(100,13).Value = -2859
$ws.Cells.Item(126,8).Value = 2833.3333
$ws.Cells.Item(126,9).Value = 3000
$ws.Cells.Item(126,11).Value = 9000
$ws.Cells.Item(126,13).Value = -6530

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(88,8).Value = 41725.668
$ws.Cells.Item(88,10).Value = 41725.668
$ws.Cells.Item(88,12).Value = 41725.668
$ws.Cells.Item(88,14).Value = -42537.668
$ws.Cells.Item(91,8).Value = 41725.668
$ws.Cells.Item(91,10).Value = 41725.668
$ws.Cells.Item(91,12).Value = 41725.668
$ws.Cells.Item(91,14).Value = -44533.668
$ws.Cells.Item(113,8).Value = 766.8182
$ws.Cells.Item(113,9).Value = 700.4
$ws.Cells.Item(113,11).Value = 2101.2
$ws.Cells.Item(113,13).Value = 68.80000000000018

